$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44839
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15500
$ws.Range("P2").Value = 1192

# Row 3
$ws.Range("D3").Value = 44868
$ws.Range("K3").Value = 18000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 18000
$ws.Range("P3").Value = 1385

# Row 4
$ws.Range("D4").Value = 44915
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 18000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 18000
$ws.Range("P4").Value = 1385

# Row 5
$ws.Range("D5").Value = 44797
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = 12500
$ws.Range("P5").Value = 962

# Row 6
$ws.Range("D6").Value = 44874
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 17000
$ws.Range("M6").Value = 17000
$ws.Range("P6").Value = 1308

# Row 8
$ws.Range("D8").Value = 44959
$ws.Range("K8").Value = 19000
$ws.Range("L8").Value = 19000
$ws.Range("M8").Value = 19000
$ws.Range("P8").Value = 1462

# Row 9
$ws.Range("D9").Value = 44804
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 12000
$ws.Range("M9").Value = 12500
$ws.Range("P9").Value = 962

# Row 10
$ws.Range("D10").Value = 44846
$ws.Range("J10").Value = 30

# Row 12
$ws.Range("D12").Value = 44943
$ws.Range("J12").Value = 30
$ws.Range("K12").Value = 17000
$ws.Range("L12").Value = 17000
$ws.Range("M12").Value = 17000
$ws.Range("P12").Value = 1308

# Row 13
$ws.Range("D13").Value = 44895

# Row 14
$ws.Range("D14").Value = 44810
$ws.Range("J14").Value = 40

# Row 15
$ws.Range("D15").Value = 45155
$ws.Range("K15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 20000
$ws.Range("P15").Value = 1538

# Row 16
$ws.Range("D16").Value = 44894
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = 18000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 18000
$ws.Range("P16").Value = 1385

# Row 17
$ws.Range("D17").Value = 44922
$ws.Range("K17").Value = 17000
$ws.Range("L17").Value = 17000
$ws.Range("M17").Value = 17000
$ws.Range("P17").Value = 1308

# Row 18
$ws.Range("D18").Value = 44841

# Row 19
$ws.Range("D19").Value = 44859
$ws.Range("K19").Value = 13000
$ws.Range("L19").Value = 13000
$ws.Range("M19").Value = 13000
$ws.Range("P19").Value = 1000

# Row 20
$ws.Range("D20").Value = 44832
$ws.Range("J20").Value = 60
$ws.Range("K20").Value = 17000
$ws.Range("M20").Value = 17500
$ws.Range("P20").Value = 1346
